# Revert "Revert "UI changes""
# The "password" and "password_confirmation" columns (F:G) are removed again,
# and the enrollment_id for the first data row (row 2 / rwilliams) is restored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the password / password_confirmation columns entirely (F:G), along
# with their placeholder 1234 values, shifting nothing else.
$ws.Range("F1:G4").EntireColumn.Delete()

# Restore the missing enrollment_id value for rwilliams.
$ws.Range("E2").Value = 1010

# Leave the selection where the original edit left it.
$ws.Range("E7").Select()
